$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")
$ws.Range("B1").Value = 1
$ws.Range("C1").Value = 1
$ws.Range("C2").Value = 1
$ws.Range("C3").Value = 1
$ws.Range("C4").Select()
$excel.CalculateFull()
